$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> new value for column F (dSF)
$updates = @{
    5  = -5
    6  = -4
    12 = -1
    14 = -3
    15 = -1
    23 = 4
    27 = -3
    28 = -1
    36 = 6
    37 = -5
    41 = -1
    49 = 2
    50 = -5
    57 = 6
    61 = 2
    62 = 2
    67 = -4
    69 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
